# Daily update at 8 AM UTC
# Appends the next day's row (row 41) to the "Wins Over Time" tracking sheet
# and moves the "last row" date-format style (YYYY-MM-DD, no time) down to
# the newly appended row, restoring the previous last row (40) to the
# standard YYYY-MM-DD HH:MM:SS date format used by all other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 was the last row and used the "short date" style; since a new
# row is being appended below it, it reverts to the regular date/time style.
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 41: next day's data.
$ws.Range("A41").Value = 45626
$ws.Range("A41").NumberFormat = "YYYY-MM-DD"

$ws.Range("B41").Value = 108
$ws.Range("C41").Value = 89
$ws.Range("D41").Value = 96
